# FoodTracker.xlsx — finalize Friday (row 29) dinner entry + totals,
# and update the earlier Thursday (row 28) "macrosTotal" figure.
# Also refresh the sheet's saved view position/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 28 (Thursday): macrosTotal (AI28) recorded value changes ---
$ws.Range("AI28").Value2 = "236.6,153.4,82.4"

# --- Row 29 (Friday): add the Dinner entry and fill out the day's totals ---
$ws.Range("I29").Value2 = "(H)Pigeon peas Dal with Protien tortillas, peanuts"
$ws.Range("J29").Value2 = "90.0,54,35.8"

$ws.Range("AA29").Value2 = 264.7
$ws.Range("AB29").Value2 = 150.4
$ws.Range("AC29").Value2 = 91.9
$ws.Range("AD29").Value2 = 2059
$ws.Range("AE29").Value2 = 1089
$ws.Range("AF29").Formula = "=AD29-AE29"
$ws.Range("AG29").Value2 = 3
$ws.Range("AH29").Value2 = 16719
$ws.Range("AI29").Value2 = "264.7,150.4,91.9"

# --- Restore the saved scroll position / selection ---
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 24
$win.ScrollRow = 4
$ws.Range("AG29").Select()

$wb.Saved = $false
